$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: several columns in this sheet store numeric-looking values
# (order-limit / sell-price) as plain TEXT inside cells that nonetheless
# keep a numeric display format. Assigning a numeric-looking string via
# .Value normally makes the engine coerce it to a real number (and
# re-point the cell at a plain "text" style in the process), which loses
# both the original style id and the "t=s" shared-string typing. Round-
# tripping the NumberFormat forces the assignment to stick as literal
# text while leaving the cell's number format (and therefore style id)
# unchanged.
function Set-TextValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# --- Step 1: insert a new row before the existing row 7 (DEPAKINE),
#     copying row 7's formatting (styles + merges) so the new row matches
#     the visual/style layout already used by the table body rows.
$ws.Range("A7:Q7").Insert()
$ws.Range("A8:Q8").Copy($ws.Range("A7:Q7"))

# Step 1b: fill in the new row 7 (COGINTOL 20 TAB.)
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "COGINTOL 20 TAB."
$ws.Range("H7").Value = "0:0"
Set-TextValue $ws.Range("L7") "1"
$ws.Range("N7").Value = "40.00"
Set-TextValue $ws.Range("P7") "40.0000"
$ws.Range("Q7").Value = "1:0"

# Row 8 is now the original DEPAKINE row (shifted down). Update its stock
# value + index number; other values stay the same.
$ws.Range("A8").Value = 2
$ws.Range("H8").Value = "0:0"

# Row 9 is now the original ERASTAPEX row (shifted down). Its index
# number moves from 2 to 3 and its stock value changes 1:2 -> 0:2.
$ws.Range("A9").Value = 3
$ws.Range("H9").Value = "0:2"

# --- Step 2: insert a new row after row 9 (ERASTAPEX) and before the
#     total row, copying row 9's formatting.
$ws.Range("A10:Q10").Insert()
$ws.Range("A9:Q9").Copy($ws.Range("A10:Q10"))

# Step 2b: fill in the new row 10 (HALONACE 5 MG 10 TABS.)
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "HALONACE 5 MG 10 TABS."
$ws.Range("H10").Value = "0:0"
Set-TextValue $ws.Range("L10") "1"
$ws.Range("N10").Value = "17.00"
Set-TextValue $ws.Range("P10") "17.0000"
$ws.Range("Q10").Value = "1:0"

# --- Step 3: the total row (was row 9, now row 11) and the footer row
#     (was row 10, now row 12) need updated values.
$ws.Range("P11").Value = 315
$ws.Range("A12").Value = "Monday, 15 September, 2025 9:36 AM"

# --- Step 4: restore the expected row heights for the table body rows.
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 16.5

Write-Output "done"
